{"js": "// 1. Remove the stray \"_GoBack\" bookmark left over from the previous save.\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\n\n// 2. Merge the three runs (split apart by grammar-check proofErr markers)\n//    that make up the \">>> your stuff after this line >>>\" placeholder\n//    paragraph back into a single run, and insert the new \"Version\n//    management\" paragraph right after it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst placeholder = paragraphs.items.find((p) => p.text.indexOf(\">>>\") !== -1 && p.text.indexOf(\"your stuff after this line\") !== -1);\n\nplaceholder.getRange().insertText(\">>>  your stuff after this line >>>\", Word.InsertLocation.replace);\n\nplaceholder.insertParagraph(\n  \"Version management help the developer to see all the changes that happen to the file last time it was edit, what was added what was removed. This should be very helpful since unwanted changes can be undo and good changes can be kept. Its also a form of backup for depends on what you want to do.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark left over from the previous save.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Merge the three runs (split apart by grammar-check proofErr markers)\n#    that make up the \">>> your stuff after this line >>>\" placeholder\n#    paragraph back into a single run (no more proofErr splits), then\n#    insert the new \"Version management\" paragraph right after it.\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text -like \"*your stuff after this line*\") {\n\n        $find = $r.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute(\n            \">>>  your stuff after this line >>>\",\n            $false, $false, $false, $false, $false, $true, 1, $false,\n            \">>>  your stuff after this line >>>\",\n            2\n        )\n\n        $r.InsertParagraphAfter()\n        $newPara = $p.Next()\n        $newPara.Range.Text = \"Version management help the developer to see all the changes that happen to the file last time it was edit, what was added what was removed. This should be very helpful since unwanted changes can be undo and good changes can be kept. Its also a form of backup for depends on what you want to do.\"\n        break\n    }\n}\n"}
